# Update stats for 2026-02 (row 27)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B27").Value = 6549
$ws.Range("C27").Value = 1016
$ws.Range("D27").Value = 6105506
$ws.Range("E27").Value = 932.2806535348908
$ws.Range("F27").Value = 10.0672268907563
$ws.Range("G27").Value = 7.286166842661035
$ws.Range("H27").Value = 25.22812499294944
